# Apply updated crypto price/volume data as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.866.11"
$ws.Range("E2").Value = "  -1.88%  "

$ws.Range("D3").Value = "'1.803.03"
$ws.Range("E3").Value = "  -1.25%  "

$ws.Range("D4").Value = "'0.9999"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'309.20"
$ws.Range("E5").Value = "  -1.89%  "

$ws.Range("D6").Value = "'0.9999"

$ws.Range("D7").Value = "'0.4658"
$ws.Range("E7").Value = "  +4.06%  "

$ws.Range("D8").Value = "'0.3701"
$ws.Range("E8").Value = "  -2.11%  "

$ws.Range("D9").Value = "'0.07362"
$ws.Range("E9").Value = "  -1.59%  "

$ws.Range("D10").Value = "'0.8680"
$ws.Range("E10").Value = "  -2.16%  "

$ws.Range("D11").Value = "'20.37"
$ws.Range("E11").Value = "  -3.15%  "

$ws.Range("D12").Value = "'1.784.59"
$ws.Range("E12").Value = "  -2.35%  "

$ws.Range("D13").Value = "'5.359"
$ws.Range("E13").Value = "  -1.83%  "

$ws.Range("D14").Value = "'92.29"
$ws.Range("E14").Value = "  -1.73%  "

$ws.Range("D15").Value = "'6.519"
$ws.Range("E15").Value = "  -3.48%  "

$ws.Range("D16").Value = "'0.07028"
$ws.Range("E16").Value = "  -1.26%  "

$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -0.04%  "

$ws.Range("D18").Value = "'0.000008701"
$ws.Range("E18").Value = "  -1.06%  "

$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  -0.03%  "

$ws.Range("D20").Value = "'14.67"
$ws.Range("E20").Value = "  -3.21%  "

$ws.Range("D21").Value = "'26.862.14"
$ws.Range("E21").Value = "  -1.94%  "

$ws.Range("D22").Value = "'5.286"
$ws.Range("E22").Value = "  -2.30%  "

$ws.Range("D23").Value = "'10.57"
$ws.Range("E23").Value = "  -3.95%  "

$ws.Range("D24").Value = "'2.018.72"
$ws.Range("E24").Value = "  -2.03%  "

$ws.Range("D25").Value = "'1.905"
$ws.Range("E25").Value = "  -2.98%  "

$ws.Range("D26").Value = "'151.57"
$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("D27").Value = "'18.33"
$ws.Range("E27").Value = "  -2.14%  "

$ws.Range("D28").Value = "'2.141"
$ws.Range("E28").Value = "  -7.22%  "

$ws.Range("D29").Value = "'5.250"
$ws.Range("E29").Value = "  -3.12%  "

$ws.Range("D30").Value = "'116.03"
$ws.Range("E30").Value = "  -1.56%  "

$ws.Range("D31").Value = "'0.08920"
$ws.Range("E31").Value = "  +0.29%  "

$ws.Range("D32").Value = "'0.7609"
$ws.Range("E32").Value = "  -3.89%  "

$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Value = "'2.937"
$ws.Range("E33").Value = "  +0.30%  "

$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.151"
$ws.Range("E34").Value = "  -4.76%  "

$ws.Range("D35").Value = "'4.462"
$ws.Range("E35").Value = "  -3.11%  "

$ws.Range("D36").Value = "'0.9996"
$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("E37").Value = "  -0.87%  "

$ws.Range("D38").Value = "'0.01954"
$ws.Range("E38").Value = "  -1.90%  "

$ws.Range("D39").Value = "'0.05248"
$ws.Range("E39").Value = "  -1.24%  "

$ws.Range("D40").Value = "'2.929"
$ws.Range("E40").Value = "  +1.86%  "

$ws.Range("D41").Value = "'7.220"
$ws.Range("E41").Value = "  -1.41%  "

$ws.Range("D42").Value = "'0.5302"
$ws.Range("E42").Value = "  -1.00%  "

$ws.Range("D43").Value = "'2.362"
$ws.Range("E43").Value = "  +1.52%  "

$ws.Range("E44").Value = "  -3.57%  "

$ws.Range("D45").Value = "'8.509"
$ws.Range("E45").Value = "  -2.01%  "

$ws.Range("D46").Value = "'0.5015"
$ws.Range("E46").Value = "  -2.09%  "

$ws.Range("D47").Value = "'10.29"
$ws.Range("E47").Value = "  -3.59%  "

$ws.Range("D48").Value = "'104.07"
$ws.Range("E48").Value = "  -1.31%  "

$ws.Range("D49").Value = "'0.9994"
$ws.Range("E49").Value = "  -0.10%  "

$ws.Range("D50").Value = "'1.660"
$ws.Range("E50").Value = "  -2.22%  "

$ws.Range("D51").Value = "'0.06285"
$ws.Range("E51").Value = "  -1.99%  "
